$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 10:46:31"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-26 10:46:27"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45bb93908cbba444cceaab03faf865e08e1e4f62/e2e/63679998-18cd-485d-97e6-76d9b68749bd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06ebd040b900c486a3ffcf703d7b4ff012346219/e2e/63679998-18cd-485d-97e6-76d9b68749bd.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-26 10:46:31"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45bb93908cbba444cceaab03faf865e08e1e4f62/e2e/63679998-18cd-485d-97e6-76d9b68749bd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06ebd040b900c486a3ffcf703d7b4ff012346219/e2e/63679998-18cd-485d-97e6-76d9b68749bd.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
